$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

function Get-Para2Text {
    $pp = $tr.Paragraphs(2, 1)
    return $pp.Text
}

function Get-Para2Start {
    $pp = $tr.Paragraphs(2, 1)
    return $pp.Start
}

# The shared "Calibri (Body)" typeface value every run in this paragraph
# already carries explicitly. Re-assigning it to a sub-range forces the
# engine to split a run at that boundary without changing its appearance
# or introducing any new formatting attribute.
$latin = "Calibri (Body)"

# ---------------------------------------------------------------------
# Edit A: "... dari permainan board game yang ..."
#      -> "... dari suatu permainan dalam board game yang ..."
# ---------------------------------------------------------------------

# A1: insert "suatu " using the existing (non-flagged) space run that sits
#     right before "permainan" as the formatting donor, so the inserted
#     space keeps that run's formatting.
$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf(" permainan board game")
$spaceBefore = $tr.Characters($pstart + $idx, 1)
$spaceBefore.InsertAfter("suatu ")

# Split the merged "suatu permainan" run (which took on "permainan"'s own
# formatting) into separate "suatu" / " " / "permainan" runs.
$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf("suatu permainan")
$rSuatu = $tr.Characters($pstart + $idx, 5)
$rSuatu.Font.Name = $latin
$rSp1 = $tr.Characters($pstart + $idx + 5, 1)
$rSp1.Font.Name = $latin

# A2: insert " dalam" after "permainan", using the existing (non-flagged)
#     " board game yang " run as the donor for the new trailing space.
$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf(" board game yang ")
$phraseRun = $tr.Characters($pstart + $idx, 1)
$phraseRun.InsertBefore(" dalam")

# Split off the leading " " and "dalam" from the merged "permainan dalam"
# piece (which took "permainan"'s own formatting).
$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf(" dalam board game")
$rSp2 = $tr.Characters($pstart + $idx, 1)
$rSp2.Font.Name = $latin
$rDalam = $tr.Characters($pstart + $idx + 1, 5)
$rDalam.Font.Name = $latin

# ---------------------------------------------------------------------
# Edit B: "... game yang ditentukan admin ..."
#      -> "... game yang telah ditentukan admin ..."
# ---------------------------------------------------------------------
$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf(" ditentukan")
$spaceBefore2 = $tr.Characters($pstart + $idx, 1)
$spaceBefore2.InsertAfter("telah ")

$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf("telah ditentukan")
$rTelah = $tr.Characters($pstart + $idx, 5)
$rTelah.Font.Name = $latin
$rSp3 = $tr.Characters($pstart + $idx + 5, 1)
$rSp3.Font.Name = $latin

# ---------------------------------------------------------------------
# Edit C: "... dibuat untuk OS Android."
#      -> "... dibuat untuk OS berbasis Android."
# ---------------------------------------------------------------------
$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf(" OS Android.")
$rOsSpace = $tr.Characters($pstart + $idx + 3, 1)   # the space right before "Android."
$rOsSpace.InsertAfter("berbasis ")

$pstart = Get-Para2Start
$t = Get-Para2Text
$idx = $t.IndexOf("berbasis Android.")
$rBerbasis = $tr.Characters($pstart + $idx, 9)
$rBerbasis.Font.Name = $latin
